# lesson-21.xlsx — reorder the English/Japanese word-list rows.
#
# The data rows (2-48 and 58-119) were shuffled into a new order while the
# header row (1) and the "far (away) .. anyhow; anyway" block (49-57) stayed
# put. The new order is made up of contiguous chunks lifted from the old
# sheet, so we read every source chunk into memory first (several
# sources/destinations overlap) and then write them all back out in their
# new positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- read all source blocks (order doesn't matter, nothing is written yet) ---
$blockA = $ws.Range("A24:B33").Value2   # to notice .. to step on            -> A2:B11
$blockB = $ws.Range("A13:B23").Value2   # government .. absence; not at home -> A12:B22
$blockC = $ws.Range("A34:B48").Value2   # to turn down .. to contact (old)   -> A23:B37
$blockD = $ws.Range("A2:B12").Value2    # to bully .. to contact (old)       -> A38:B48

$blockE = $ws.Range("A81:B93").Value2   # to believe .. kitchen              -> A58:B70
$blockF = $ws.Range("A107:B119").Value2 # to miss (transportation) .. sister school -> A71:B83
$blockG = $ws.Range("A58:B70").Value2   # at first .. kitchen (old)          -> A84:B96
$blockH = $ws.Range("A94:B106").Value2  # happy (しあわ) .. to be late (する) -> A97:B109
$blockI = $ws.Range("A71:B80").Value2   # worried about .. sister school(old)-> A110:B119

# --- write every block to its new home ---
$ws.Range("A2:B11").Value2 = $blockA
$ws.Range("A12:B22").Value2 = $blockB
$ws.Range("A23:B37").Value2 = $blockC
$ws.Range("A38:B48").Value2 = $blockD

$ws.Range("A58:B70").Value2 = $blockE
$ws.Range("A71:B83").Value2 = $blockF
$ws.Range("A84:B96").Value2 = $blockG
$ws.Range("A97:B109").Value2 = $blockH
$ws.Range("A110:B119").Value2 = $blockI
